$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the report time-range header text
$ws.Range("A1").Value = "From: 28/01/2018 at 1745`nTo: 28/01/2018 at 1840"
$ws.Rows.Item(1).AutoFit()

# Helper to set a row of percentage values across columns A:F
function Set-Row($r, $vals) {
    for ($i = 0; $i -lt $vals.Length; $i++) {
        if ($null -ne $vals[$i]) {
            $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
        }
    }
}

Set-Row 6  @(99.67, 99.13, 99.56, 99.34, 99.29, 96.89)
Set-Row 12 @(99.67, 99.13, 99.56, 99.34, 99.56, 99.56)
Set-Row 18 @(99.67, 99.13, 99.56, 99.34, 99.56, 90.82)
Set-Row 24 @(99.67, 99.13, 99.56, 99.34, 99.56, 90.82)
Set-Row 28 @(99.4,  97.87, 97.49, 93.45, 91.97, 92.79)
Set-Row 34 @(99.67, 99.13, 99.56, 99.34, 99.56, 90.82)
Set-Row 38 @(95.85, 99.4,  $null, $null, $null, $null)
Set-Row 44 @(99.67, 99.13, 99.56, 99.34, 99.56, 96.45)
Set-Row 48 @(88.42, $null, $null, $null, $null, $null)
Set-Row 54 @(99.67, 99.13, 99.56, 99.34, 99.56, 98.85)
Set-Row 58 @(92.03, 91.81, 93.34, 92.63, $null, $null)
Set-Row 64 @(99.67, 99.13, 99.56, 99.34, 99.56, 98.85)
Set-Row 68 @(92.03, 91.81, 92.08, 92.57, $null, $null)
Set-Row 74 @(99.67, 99.13, 99.56, 99.34, 99.56, 98.85)
Set-Row 78 @(92.03, 91.81, 93.77, 86.07, 85.91, 84.76)
Set-Row 84 @(99.67, 99.13, 99.56, 99.34, 99.56, 98.85)
Set-Row 88 @(92.19, 90.01, $null, $null, $null, $null)
Set-Row 94 @(99.67, 99.13, 99.56, 99.34, 99.56, 98.85)
Set-Row 98 @(92.03, 91.81, 93.77, 86.07, 79.63, $null)
